$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.165.27"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "1.556.93"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'206.29"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'22.15"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").Value = "'0.0592"
$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "1.778.04"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "1.555.38"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "'0.515"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").Value = "'62.77"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "27.143.47"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "'214.90"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "'9.33"
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("D24").Value = "'1.99"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'151.90"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "'14.88"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "'0.0461"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").Value = "1.377.11"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "'0.949"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "'0.810"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").Value = "'0.513"
$ws.Range("E40").Value = "  -4.52%  "
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("E43").Value = "  +3.72%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'63.15"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "1.690.83"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").Value = "'85.28"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  +0.10%  "
